# carbatpy/heat_exchanger_input0.xlsx - Excel parameter layout rework
# Reorganises the "variable name"/"comment/choice" columns into a
# consistent variable_name/value/value_min/value_max/fixed/name_fluid/unit/
# option/comment/value_choices layout across all four sheets and adds a
# "configuration" (Geometry) / extra rows (Problem description) at the end.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Fluid 1
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Fluid 1")
$ws2.Activate()
$ws2.UsedRange.Clear()

$ws2.Cells.Item(1,1).Value = "variable_name"
$ws2.Cells.Item(1,2).Value = "value"
$ws2.Cells.Item(1,3).Value = "value_min"
$ws2.Cells.Item(1,4).Value = "value_max"
$ws2.Cells.Item(1,5).Value = "fixed"
$ws2.Cells.Item(1,6).Value = "name_fluid"
$ws2.Cells.Item(1,7).Value = "unit"
$ws2.Cells.Item(1,8).Value = "option"
$ws2.Cells.Item(1,9).Value = "comment"
$ws2.Cells.Item(1,10).Value = "value_choices"

$ws2.Cells.Item(2,1).Value = "props"
$ws2.Cells.Item(2,2).Value = "REFPROP"

$ws2.Cells.Item(3,1).Value = "number_compounds"
$ws2.Cells.Item(3,2).Value = 2

$ws2.Cells.Item(4,1).Value = "fl1"
$ws2.Cells.Item(4,2).Value = 0.6
$ws2.Cells.Item(4,3).Value = 0.1
$ws2.Cells.Item(4,4).Value = 0.8
$ws2.Cells.Item(4,5).Value = "'True"
$ws2.Cells.Item(4,6).Value = "Propane"
$ws2.Cells.Item(4,7).Value = "mole fraction"

$ws2.Cells.Item(5,1).Value = "fl2"
$ws2.Cells.Item(5,2).Value = 0.4
$ws2.Cells.Item(5,3).Value = 0.2
$ws2.Cells.Item(5,4).Value = 0.9
$ws2.Cells.Item(5,5).Value = "'True"
$ws2.Cells.Item(5,6).Value = "Pentane"
$ws2.Cells.Item(5,7).Value = "mole fraction"

$ws2.Cells.Item(6,1).Value = "T_in"
$ws2.Cells.Item(6,2).Value = 350
$ws2.Cells.Item(6,3).Value = 290
$ws2.Cells.Item(6,4).Value = 370
$ws2.Cells.Item(6,5).Value = "'True"
$ws2.Cells.Item(6,7).Value = "K"

$ws2.Cells.Item(7,1).Value = "p_in"
$ws2.Cells.Item(7,2).Value = 1000000
$ws2.Cells.Item(7,2).NumberFormat = "0.00E+00"
$ws2.Cells.Item(7,3).Value = 100000
$ws2.Cells.Item(7,3).NumberFormat = "0.00E+00"
$ws2.Cells.Item(7,4).Value = 2000000
$ws2.Cells.Item(7,4).NumberFormat = "0.00E+00"
$ws2.Cells.Item(7,5).Value = "'True"
$ws2.Cells.Item(7,7).Value = "Pa"

$ws2.Cells.Item(8,1).Value = "m_dot"
$ws2.Cells.Item(8,2).Value = 0.012
$ws2.Cells.Item(8,3).Value = 0.01
$ws2.Cells.Item(8,4).Value = 0.02
$ws2.Cells.Item(8,7).Value = "kg/s"

$ws2.Columns.Item(1).ColumnWidth = 16

$ws2.Range("G1:H1048576").Select()

# ---------------------------------------------------------------------
# Sheet: Fluid 2
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Fluid 2")
$ws3.Activate()
$ws3.UsedRange.Clear()

$ws3.Cells.Item(1,1).Value = "variable_name"
$ws3.Cells.Item(1,2).Value = "value"
$ws3.Cells.Item(1,3).Value = "value_min"
$ws3.Cells.Item(1,4).Value = "value_max"
$ws3.Cells.Item(1,5).Value = "fixed"
$ws3.Cells.Item(1,6).Value = "name_fluid"
$ws3.Cells.Item(1,7).Value = "unit"
$ws3.Cells.Item(1,8).Value = "option"
$ws3.Cells.Item(1,9).Value = "comment"
$ws3.Cells.Item(1,10).Value = "value_choices"

$ws3.Cells.Item(2,1).Value = "props"
$ws3.Cells.Item(2,2).Value = "REFPROP"

$ws3.Cells.Item(3,1).Value = "number_compounds"
$ws3.Cells.Item(3,2).Value = 1

$ws3.Cells.Item(4,1).Value = "fl1"
$ws3.Cells.Item(4,2).Value = 1
$ws3.Cells.Item(4,3).Value = 1
$ws3.Cells.Item(4,4).Value = 1
$ws3.Cells.Item(4,5).Value = "'True"
$ws3.Cells.Item(4,6).Value = "Water"
$ws3.Cells.Item(4,7).Value = "mole fraction"

$ws3.Cells.Item(5,1).Value = "T_in"
$ws3.Cells.Item(5,2).Value = 350
$ws3.Cells.Item(5,3).Value = 290
$ws3.Cells.Item(5,4).Value = 370
$ws3.Cells.Item(5,5).Value = "'True"
$ws3.Cells.Item(5,7).Value = "mole fraction"

$ws3.Cells.Item(6,1).Value = "p_in"
$ws3.Cells.Item(6,2).Value = 1000000
$ws3.Cells.Item(6,2).NumberFormat = "0.00E+00"
$ws3.Cells.Item(6,3).Value = 100000
$ws3.Cells.Item(6,3).NumberFormat = "0.00E+00"
$ws3.Cells.Item(6,4).Value = 2000000
$ws3.Cells.Item(6,4).NumberFormat = "0.00E+00"
$ws3.Cells.Item(6,5).Value = "'True"
$ws3.Cells.Item(6,7).Value = "K"

$ws3.Cells.Item(7,1).Value = "m_dot"
$ws3.Cells.Item(7,2).Value = 0.012
$ws3.Cells.Item(7,3).Value = 0.01
$ws3.Cells.Item(7,4).Value = 0.02
$ws3.Cells.Item(7,7).Value = "Pa"

$ws3.Cells.Item(8,7).Value = "kg/s"

# empty, style-only cells that were already present in the source file
$ws3.Cells.Item(9,2).NumberFormat = "0.00E+00"
$ws3.Cells.Item(12,2).NumberFormat = "0.00E+00"
$ws3.Cells.Item(13,2).NumberFormat = "0.00E+00"

$ws3.Range("I8").Select()

# ---------------------------------------------------------------------
# Sheet: Problem description
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Problem description")
$ws4.Activate()
$ws4.UsedRange.Clear()

$ws4.Cells.Item(1,1).Value = "variable_name"
$ws4.Cells.Item(1,2).Value = "value"
$ws4.Cells.Item(1,3).Value = "value_min"
$ws4.Cells.Item(1,4).Value = "value_max"
$ws4.Cells.Item(1,5).Value = "fixed"
$ws4.Cells.Item(1,6).Value = "unused"
$ws4.Cells.Item(1,7).Value = "option"
$ws4.Cells.Item(1,8).Value = "unit"
$ws4.Cells.Item(1,9).Value = "comment"
$ws4.Cells.Item(1,10).Value = "value_choices"

$ws4.Cells.Item(2,1).Value = "calculation_type"
$ws4.Cells.Item(2,2).Value = "output"

$ws4.Cells.Item(3,1).Value = "T_ref"
$ws4.Cells.Item(3,2).Value = 283.15
$ws4.Cells.Item(3,8).Value = "K"

$ws4.Range("F2").Select()

# ---------------------------------------------------------------------
# Sheet: Geometry  (kept active / tab-selected, like in the source file)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Geometry")
$ws1.Activate()
$ws1.UsedRange.Clear()

$ws1.Cells.Item(1,1).Value = "variable_name"
$ws1.Cells.Item(1,2).Value = "value_1"
$ws1.Cells.Item(1,3).Value = "value_2"
$ws1.Cells.Item(1,4).Value = "fixed_1"
$ws1.Cells.Item(1,5).Value = "fixed_2"
$ws1.Cells.Item(1,6).Value = "unit"
$ws1.Cells.Item(1,7).Value = "option"
$ws1.Cells.Item(1,8).Value = "comment"

$ws1.Cells.Item(2,1).Value = "hex_type"
$ws1.Cells.Item(2,2).Value = "counterflow"

$ws1.Cells.Item(3,1).Value = "specification"
$ws1.Cells.Item(3,2).Value = "shellTube"

$ws1.Cells.Item(4,1).Value = "U"
$ws1.Cells.Item(4,2).Value = 500
$ws1.Cells.Item(4,3).Value = 2
$ws1.Cells.Item(4,4).Value = "'True"
$ws1.Cells.Item(4,5).Value = "'True"
$ws1.Cells.Item(4,6).Value = "W/ m2 / K"

$ws1.Cells.Item(5,1).Value = "tubes"
$ws1.Cells.Item(5,2).Value = 12
$ws1.Cells.Item(5,3).Value = 1
$ws1.Cells.Item(5,4).Value = "'True"
$ws1.Cells.Item(5,5).Value = "'True"

$ws1.Cells.Item(6,1).Value = "d_in"
$ws1.Cells.Item(6,2).Value = 0.01
$ws1.Cells.Item(6,2).NumberFormat = "0.00E+00"
$ws1.Cells.Item(6,3).Value = 0.08
$ws1.Cells.Item(6,3).NumberFormat = "0.00E+00"
$ws1.Cells.Item(6,4).Value = "'True"
$ws1.Cells.Item(6,5).Value = "'True"
$ws1.Cells.Item(6,6).Value = "m"

$ws1.Cells.Item(7,1).Value = "length"
$ws1.Cells.Item(7,2).Value = 4
$ws1.Cells.Item(7,4).Value = "'True"
$ws1.Cells.Item(7,5).Value = "'True"
$ws1.Cells.Item(7,6).Value = "m"

$ws1.Cells.Item(8,1).Value = "wall_thickness"
$ws1.Cells.Item(8,2).Value = 0.001
$ws1.Cells.Item(8,2).NumberFormat = "0.00E+00"
$ws1.Cells.Item(8,3).Value = 0.005
$ws1.Cells.Item(8,3).NumberFormat = "0.00E+00"
$ws1.Cells.Item(8,4).Value = "'True"
$ws1.Cells.Item(8,5).Value = "'True"
$ws1.Cells.Item(8,6).Value = "m"

$ws1.Cells.Item(9,1).Value = "cp"
$ws1.Cells.Item(9,2).Value = 800
$ws1.Cells.Item(9,3).Value = 800
$ws1.Cells.Item(9,4).Value = "'True"
$ws1.Cells.Item(9,5).Value = "'True"
$ws1.Cells.Item(9,6).Value = "J/kg/K"

$ws1.Cells.Item(10,1).Value = "rho"
$ws1.Cells.Item(10,2).Value = 8000
$ws1.Cells.Item(10,2).NumberFormat = "0.00E+00"
$ws1.Cells.Item(10,3).Value = 8000
$ws1.Cells.Item(10,3).NumberFormat = "0.00E+00"
$ws1.Cells.Item(10,4).Value = "'True"
$ws1.Cells.Item(10,5).Value = "'True"
$ws1.Cells.Item(10,6).Value = "kg/m3"

$ws1.Cells.Item(11,1).Value = "material"
$ws1.Cells.Item(11,2).Value = "Steel"
$ws1.Cells.Item(11,3).Value = "Steel"
$ws1.Cells.Item(11,4).Value = "'True"
$ws1.Cells.Item(11,5).Value = "'True"
$ws1.Cells.Item(11,6).Value = "W / m / K"

$ws1.Cells.Item(12,1).Value = "thermal_conductivity"
$ws1.Cells.Item(12,2).Value = 16
$ws1.Cells.Item(12,2).NumberFormat = "0.00E+00"
$ws1.Cells.Item(12,3).Value = 16
$ws1.Cells.Item(12,3).NumberFormat = "0.00E+00"
$ws1.Cells.Item(12,4).Value = "'True"
$ws1.Cells.Item(12,5).Value = "'True"

$ws1.Cells.Item(13,1).Value = "configuration"
$ws1.Cells.Item(13,2).Value = "Fluid_1"
$ws1.Cells.Item(13,3).Value = "constant"
$ws1.Cells.Item(13,8).Value = "Name of Sheet for the selection of 1 or two fluids,or  value_2: constant or T-slope"

# Freeze header row + first column, matching the recorded view state
$ws1.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws1.Range("H13").Select()

$ws1.Activate()
